$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 217.5
$ws.Range("I9").Value = 287.625
$ws.Range("J9").Value = 77.25
$ws.Range("K9").Value = 287.625
$ws.Range("L9").Value = 77.25
$ws.Range("M9").Value = -118.625
$ws.Range("N9").Value = -415.25
$ws.Range("H17").Value = 170333.33
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 170333.33
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 510999.99
$ws.Range("N17").Value = -511335.99
$ws.Range("H28").Value = 2824
$ws.Range("I28").Value = 432
$ws.Range("J28").Value = 10000
$ws.Range("K28").Value = 432
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = 53
$ws.Range("N28").Value = -10970
$ws.Range("H34").Value = 12166.5
$ws.Range("I34").Value = 6875
$ws.Range("J34").Value = 22749.5
$ws.Range("K34").Value = 6875
$ws.Range("L34").Value = 22749.5
$ws.Range("M34").Value = -6672
$ws.Range("N34").Value = -23155.5
$ws.Range("H36").Value = 12166.5
$ws.Range("I36").Value = 6875
$ws.Range("J36").Value = 22749.5
$ws.Range("K36").Value = 6875
$ws.Range("L36").Value = 22749.5
$ws.Range("M36").Value = -6160
$ws.Range("N36").Value = -24179.5
$ws.Range("H70").Value = 144942.72
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 144942.72
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 434828.16
$ws.Range("N70").Value = -435368.16
$ws.Range("H73").Value = 144942.72
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 144942.72
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 434828.16
$ws.Range("N73").Value = -436700.16
$ws.Range("H76").Value = 15672.429
$ws.Range("I76").Value = 13284.5
$ws.Range("J76").Value = 30000
$ws.Range("K76").Value = 13284.5
$ws.Range("L76").Value = 30000
$ws.Range("M76").Value = -12969.5
$ws.Range("N76").Value = -30630
$ws.Range("H79").Value = 15672.429
$ws.Range("I79").Value = 13284.5
$ws.Range("J79").Value = 30000
$ws.Range("K79").Value = 13284.5
$ws.Range("L79").Value = 30000
$ws.Range("M79").Value = -12192.5
$ws.Range("N79").Value = -32184
$ws.Range("H107").Value = 655.1429000000001
$ws.Range("I107").Value = 747.6667
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 747.6667
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 1172.3333
$ws.Range("N107").Value = -3940
$ws.Range("H112").Value = 5055.1924
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 5055.1924
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 15165.5772
$ws.Range("N112").Value = -17381.5772
$ws.Range("H137").Value = 4012.8
$ws.Range("I137").Value = 3982.8
$ws.Range("J137").Value = 4042.8
$ws.Range("K137").Value = 11948.4
$ws.Range("L137").Value = 12128.4
$ws.Range("M137").Value = -9398.400000000001
$ws.Range("N137").Value = -17228.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 73692.08
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 73692.08
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 73692.08
$ws.Range("N98").Value = -79682.08
$ws.Range("H132").Value = 2219.4
$ws.Range("I132").Value = 2085
$ws.Range("J132").Value = 3429
$ws.Range("K132").Value = 6255
$ws.Range("L132").Value = 10287
$ws.Range("M132").Value = -3725
$ws.Range("N132").Value = -15347

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 28573104
$ws.Range("I86").Value = 50001480
$ws.Range("J86").Value = 1932.3334
$ws.Range("K86").Value = 50001480
$ws.Range("L86").Value = 1932.3334
$ws.Range("M86").Value = -50000357
$ws.Range("N86").Value = -4178.3334
$ws.Range("H89").Value = 28573104
$ws.Range("I89").Value = 50001480
$ws.Range("J89").Value = 1932.3334
$ws.Range("K89").Value = 250007400
$ws.Range("L89").Value = 9661.666999999999
$ws.Range("M89").Value = -250001784
$ws.Range("N89").Value = -20893.667
$ws.Range("H134").Value = 4900.5713
$ws.Range("I134").Value = 4048.3333
$ws.Range("J134").Value = 10014
$ws.Range("K134").Value = 12144.9999
$ws.Range("L134").Value = 30042
$ws.Range("M134").Value = -9609.999899999999
$ws.Range("N134").Value = -35112

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 16428.666
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 16428.666
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 16428.666
$ws.Range("N28").Value = -16918.666
$ws.Range("H31").Value = 6801.1333
$ws.Range("I31").Value = 5838.222
$ws.Range("J31").Value = 8245.5
$ws.Range("K31").Value = 5838.222
$ws.Range("L31").Value = 8245.5
$ws.Range("M31").Value = -5543.222
$ws.Range("N31").Value = -8835.5
$ws.Range("H34").Value = 6801.1333
$ws.Range("I34").Value = 5838.222
$ws.Range("J34").Value = 8245.5
$ws.Range("K34").Value = 5838.222
$ws.Range("L34").Value = 8245.5
$ws.Range("M34").Value = -5636.222
$ws.Range("N34").Value = -8649.5
$ws.Range("H62").Value = 5238.75
$ws.Range("I62").Value = 3527.5
$ws.Range("J62").Value = 6950
$ws.Range("K62").Value = 3527.5
$ws.Range("L62").Value = 6950
$ws.Range("M62").Value = -2903.5
$ws.Range("N62").Value = -8198
$ws.Range("H65").Value = 5238.75
$ws.Range("I65").Value = 3527.5
$ws.Range("J65").Value = 6950
$ws.Range("K65").Value = 17637.5
$ws.Range("L65").Value = 34750
$ws.Range("M65").Value = -14517.5
$ws.Range("N65").Value = -40990
$ws.Range("J107").Value = 2931.3333
$ws.Range("K107").Value = 29412312
$ws.Range("L107").Value = 2931.3333
$ws.Range("M107").Value = -29410392
$ws.Range("N107").Value = -6771.3333
$ws.Range("H134").Value = 2198.8333
$ws.Range("I134").Value = 2733.3333
$ws.Range("J134").Value = 1664.3334
$ws.Range("K134").Value = 8199.999899999999
$ws.Range("L134").Value = 4993.0002
$ws.Range("M134").Value = -5664.999899999999
$ws.Range("N134").Value = -10063.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 706.4
$ws.Range("I5").Value = 706.4
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2119.2
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2007.2
$ws.Range("N5").ClearContents()
$ws.Range("H33").Value = 294
$ws.Range("I33").Value = 162
$ws.Range("J33").Value = 426
$ws.Range("K33").Value = 972
$ws.Range("L33").Value = 2556
$ws.Range("M33").Value = -689
$ws.Range("N33").Value = -3122
$ws.Range("H68").Value = 933.3333
$ws.Range("I68").Value = 650
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 1950
$ws.Range("L68").Value = 4500
$ws.Range("M68").Value = -1139
$ws.Range("N68").Value = -6122
$ws.Range("H71").Value = 933.3333
$ws.Range("I71").Value = 650
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 5850
$ws.Range("L71").Value = 13500
$ws.Range("M71").Value = -1794
$ws.Range("N71").Value = -21612
$ws.Range("H107").Value = 1122.8572
$ws.Range("I107").Value = 665
$ws.Range("J107").Value = 1377.2222
$ws.Range("K107").Value = 1995
$ws.Range("L107").Value = 4131.6666
$ws.Range("M107").Value = -75
$ws.Range("N107").Value = -7971.6666
$ws.Range("H127").Value = 2000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 2000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 6000
$ws.Range("N127").Value = -15920
$ws.Range("H132").Value = 699.6667
$ws.Range("I132").Value = 707.5
$ws.Range("J132").Value = 684
$ws.Range("K132").Value = 6367.5
$ws.Range("L132").Value = 6156
$ws.Range("M132").Value = -3837.5
$ws.Range("N132").Value = -11216
$ws.Range("H135").Value = 706.4
$ws.Range("I135").Value = 706.4
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6357.599999999999
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3822.599999999999
$ws.Range("N135").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 45000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 45000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 45000
$ws.Range("N39").Value = -46064
$ws.Range("H126").Value = 2161.1538
$ws.Range("I126").Value = 2585.7144
$ws.Range("J126").Value = 1665.8334
$ws.Range("K126").Value = 7757.1432
$ws.Range("L126").Value = 4997.5002
$ws.Range("M126").Value = -5287.1432
$ws.Range("N126").Value = -9937.5002
$ws.Range("H132").Value = 4631.1924
$ws.Range("I132").Value = 4908.9165
$ws.Range("J132").Value = 1298.5
$ws.Range("K132").Value = 14726.7495
$ws.Range("L132").Value = 3895.5
$ws.Range("M132").Value = -12196.7495
$ws.Range("N132").Value = -8955.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 82021.89
$ws.Range("I74").Value = 18197
$ws.Range("J74").Value = 90000
$ws.Range("K74").Value = 18197
$ws.Range("L74").Value = 90000
$ws.Range("M74").Value = -17199
$ws.Range("N74").Value = -91996
$ws.Range("H77").Value = 82021.89
$ws.Range("I77").Value = 18197
$ws.Range("J77").Value = 90000
$ws.Range("K77").Value = 54591
$ws.Range("L77").Value = 270000
$ws.Range("M77").Value = -49599
$ws.Range("N77").Value = -279984
$ws.Range("H132").Value = 3250.7646
$ws.Range("I132").Value = 2434.9285
$ws.Range("J132").Value = 4243.9565
$ws.Range("K132").Value = 7304.7855
$ws.Range("L132").Value = 12731.8695
$ws.Range("M132").Value = -4774.7855
$ws.Range("N132").Value = -17791.8695
$ws.Range("H136").Value = 4245
$ws.Range("I136").Value = 4023.75
$ws.Range("J136").Value = 7342.5
$ws.Range("K136").Value = 12071.25
$ws.Range("L136").Value = 22027.5
$ws.Range("M136").Value = -9521.25
$ws.Range("N136").Value = -27127.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 8000
$ws.Range("I24").Value = 8000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 8000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -7770
$ws.Range("H29").Value = 8005
$ws.Range("I29").Value = 1010
$ws.Range("J29").Value = 15000
$ws.Range("K29").Value = 1010
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -720
$ws.Range("N29").Value = -15580
$ws.Range("H101").Value = 36999.8
$ws.Range("I101").Value = 62000
$ws.Range("J101").Value = 30749.75
$ws.Range("K101").Value = 62000
$ws.Range("L101").Value = 30749.75
$ws.Range("M101").Value = -58755
$ws.Range("N101").Value = -37239.75
$ws.Range("H122").Value = 3269.7083
$ws.Range("I122").Value = 2446.3157
$ws.Range("J122").Value = 6398.6
$ws.Range("K122").Value = 7338.9471
$ws.Range("L122").Value = 19195.8
$ws.Range("M122").Value = -4888.9471
$ws.Range("N122").Value = -24095.8
$ws.Range("H136").Value = 2569.843
$ws.Range("I136").Value = 2224.7446
$ws.Range("J136").Value = 6624.75
$ws.Range("K136").Value = 6674.2338
$ws.Range("L136").Value = 19874.25
$ws.Range("M136").Value = -4124.2338
$ws.Range("N136").Value = -24974.25
